$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 340.9
$ws.Range("I5").Value = 295.44446
$ws.Range("K5").Value = 295.44446
$ws.Range("M5").Value = -180.44446
$ws.Range("H33").Value = 97.125
$ws.Range("I33").Value = 72
$ws.Range("J33").Value = 116.666664
$ws.Range("K33").Value = 72
$ws.Range("L33").Value = 116.666664
$ws.Range("M33").Value = 157
$ws.Range("N33").Value = -574.666664
$ws.Range("H51").Value = 7000
$ws.Range("J51").Value = 7000
$ws.Range("L51").Value = 7000
$ws.Range("N51").Value = -7968
$ws.Range("H64").Value = 2984.7
$ws.Range("I64").Value = 2900.3333
$ws.Range("K64").Value = 2900.3333
$ws.Range("M64").Value = -2652.3333
$ws.Range("H67").Value = 2984.7
$ws.Range("I67").Value = 2900.3333
$ws.Range("K67").Value = 2900.3333
$ws.Range("M67").Value = -2042.3333
$ws.Range("H97").Value = 1099.5
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 2997
$ws.Range("L97").Value = 3600
$ws.Range("M97").Value = -2501
$ws.Range("N97").Value = -4592
$ws.Range("H125").Value = 1204.1666
$ws.Range("I125").Value = 1275
$ws.Range("J125").Value = 1062.5
$ws.Range("K125").Value = 11475
$ws.Range("L125").Value = 9562.5
$ws.Range("M125").Value = -9015
$ws.Range("N125").Value = -14482.5
$ws.Range("H138").Value = 2440.1628
$ws.Range("I138").Value = 2449.2693
$ws.Range("J138").Value = 2426.2354
$ws.Range("K138").Value = 7347.8079
$ws.Range("L138").Value = 7278.706200000001
$ws.Range("M138").Value = -2207.8079
$ws.Range("N138").Value = -17558.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5815031
$ws.Range("I2").Value = 5815031
$ws.Range("K2").Value = 5815031
$ws.Range("M2").Value = -5814918
$ws.Range("H11").Value = 3357334.8
$ws.Range("J11").Value = 36002
$ws.Range("L11").Value = 36002
$ws.Range("N11").Value = -36290
$ws.Range("H32").Value = 3039.7537
$ws.Range("I32").Value = 2118.9075
$ws.Range("K32").Value = 2118.9075
$ws.Range("M32").Value = -1831.9075
$ws.Range("H116").Value = 5815031
$ws.Range("I116").Value = 5815031
$ws.Range("K116").Value = 5815031
$ws.Range("M116").Value = -5812737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5815031
$ws.Range("I3").Value = 5815031
$ws.Range("K3").Value = 5815031
$ws.Range("M3").Value = -5814917
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225
$ws.Range("I31").Value = 1799.875
$ws.Range("K31").Value = 1799.875
$ws.Range("M31").Value = -1504.875
$ws.Range("H34").Value = 2225
$ws.Range("I34").Value = 1799.875
$ws.Range("K34").Value = 1799.875
$ws.Range("M34").Value = -1597.875
$ws.Range("H62").Value = 2542.7144
$ws.Range("I62").Value = 2542.7144
$ws.Range("K62").Value = 2542.7144
$ws.Range("M62").Value = -1918.7144
$ws.Range("H65").Value = 2542.7144
$ws.Range("I65").Value = 2542.7144
$ws.Range("K65").Value = 12713.572
$ws.Range("M65").Value = -9593.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 887.5714
$ws.Range("I32").Value = 799.3333
$ws.Range("J32").Value = 953.75
$ws.Range("K32").Value = 2397.9999
$ws.Range("L32").Value = 2861.25
$ws.Range("M32").Value = -2114.9999
$ws.Range("N32").Value = -3427.25
$ws.Range("H33").Value = 104.2
$ws.Range("I33").Value = 67
$ws.Range("J33").Value = 129
$ws.Range("K33").Value = 402
$ws.Range("L33").Value = 774
$ws.Range("M33").Value = -119
$ws.Range("N33").Value = -1340
$ws.Range("H104").Value = 2965.9583
$ws.Range("I104").Value = 1172.25
$ws.Range("K104").Value = 3516.75
$ws.Range("M104").Value = -895.75
$ws.Range("H105").Value = 2663.6365
$ws.Range("J105").Value = 2748.1904
$ws.Range("L105").Value = 8244.5712
$ws.Range("N105").Value = -13486.5712
$ws.Range("H129").Value = 46162.812
$ws.Range("J129").Value = 52672.57
$ws.Range("L129").Value = 158017.71
$ws.Range("N129").Value = -168017.71
$ws.Range("H131").Value = 796.71
$ws.Range("J131").Value = 814.29785
$ws.Range("L131").Value = 2442.89355
$ws.Range("N131").Value = -12522.89355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 63.833332
$ws.Range("I2").Value = 16.333334
$ws.Range("K2").Value = 16.333334
$ws.Range("M2").Value = 96.66666599999999
$ws.Range("H29").Value = 43703.43
$ws.Range("I29").Value = 23975
$ws.Range("K29").Value = 23975
$ws.Range("M29").Value = -23685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H68").Value = 3532.4285
$ws.Range("I68").Value = 3287.8333
$ws.Range("K68").Value = 3287.8333
$ws.Range("M68").Value = -2538.8333
$ws.Range("H71").Value = 3532.4285
$ws.Range("I71").Value = 3287.8333
$ws.Range("K71").Value = 16439.1665
$ws.Range("M71").Value = -12695.1665
$ws.Range("H82").Value = 2440.4285
$ws.Range("I82").Value = 1275
$ws.Range("K82").Value = 1275
$ws.Range("M82").Value = -914
$ws.Range("H85").Value = 2440.4285
$ws.Range("I85").Value = 1275
$ws.Range("K85").Value = 1275
$ws.Range("M85").Value = -27

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2967.6667
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 3501.5
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 3501.5
$ws.Range("M62").Value = -1276
$ws.Range("N62").Value = -4749.5
$ws.Range("H65").Value = 2967.6667
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 3501.5
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 3501.5
$ws.Range("M65").Value = -6380
$ws.Range("N65").Value = -23747.5
$ws.Range("H139").Value = 63560
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 66950
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 66950
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -77230
